# Weekly fruit/vegetable update: a new week of data (row 177) is inserted
# above the existing "Berenjena" records, pushing all subsequent rows
# (old 177..214) down by one (new 178..215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 177; everything below shifts down.
$ws.Rows("177:177").Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(177, 1).Value()  = 6
$ws.Cells.Item(177, 2).Value()  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(177, 3).Value()  = "Metropolitana"
$ws.Cells.Item(177, 4).Value()  = 44711
$ws.Cells.Item(177, 5).Value()  = 13
$ws.Cells.Item(177, 6).Value()  = 100112001
$ws.Cells.Item(177, 7).Value()  = "Berenjena"
$ws.Cells.Item(177, 8).Value()  = "Sin especificar"
$ws.Cells.Item(177, 9).Value()  = "Primera"
$ws.Cells.Item(177, 10).Value() = 250
$ws.Cells.Item(177, 11).Value() = 4000
$ws.Cells.Item(177, 12).Value() = 5000
$ws.Cells.Item(177, 13).Value() = 4600
$ws.Cells.Item(177, 14).Value() = "$/caja 50 unidades"
$ws.Cells.Item(177, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(177, 16).Value() = 92
$ws.Cells.Item(177, 17).Value() = 50
$ws.Cells.Item(177, 18).Value() = "Hortaliza"
